$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.199.69"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.564.94"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9992"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.06"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3816"
$ws.Range("E7").Value = "  +3.51%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3312"
$ws.Range("E8").Value = "  +0.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.59"
$ws.Range("E9").Value = "  -6.84%  "
$ws.Range("E10").Value = "  +1.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07403"
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9998"
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.25"
$ws.Range("E13").Value = "  -1.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.858"
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.897"
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.558.80"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001098"
$ws.Range("E17").Value = "  -2.09%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06648"
$ws.Range("E18").Value = "  -1.66%  "
$ws.Range("B19").Value = "Litecoin"
$ws.Range("C19").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "85.98"
$ws.Range("E19").Value = "  -2.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.451"
$ws.Range("E20").Value = "  +1.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9987"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.18"
$ws.Range("E22").Value = "  -1.75%  "
$ws.Range("E23").Value = "  -1.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.202.46"
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.293"
$ws.Range("E25").Value = "  -3.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.537"
$ws.Range("E26").Value = "  -1.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "149.64"
$ws.Range("E27").Value = "  -3.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.20"
$ws.Range("E28").Value = "  -2.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.901"
$ws.Range("E29").Value = "  -2.14%  "
$ws.Range("B30").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C30").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.733.73"
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.08"
$ws.Range("E31").Value = "  -2.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.104"
$ws.Range("E32").Value = "  +3.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.983"
$ws.Range("E33").Value = "  -1.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.890"
$ws.Range("E34").Value = "  -5.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08269"
$ws.Range("E35").Value = "  -1.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.371"
$ws.Range("E36").Value = "  -3.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02340"
$ws.Range("E37").Value = "  -5.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.314"
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06260"
$ws.Range("E39").Value = "  -2.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2171"
$ws.Range("E40").Value = "  -3.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.251"
$ws.Range("E41").Value = "  -2.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.13"
$ws.Range("E42").Value = "  -1.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6047"
$ws.Range("E43").Value = "  -3.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9989"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.82"
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.745"
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5850"
$ws.Range("E47").Value = "  -4.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.003"
$ws.Range("E48").Value = "  -2.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.12"
$ws.Range("E49").Value = "  -2.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.178"
$ws.Range("E50").Value = "  -2.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07045"
$ws.Range("E51").Value = "  -2.27%  "
